# Time Delay Tutorial Renewal - cleanup of stray local-path image
# descriptions left behind by the PolarisOffice import on slide 2.
#
# Two pictures on slide 2 ("Picture " / id 40 and "그림 24" / id 46) carry
# an AlternativeText (OOXML <p:cNvPr descr="...">) that is just a leaked
# local temp-file path from the authoring machine
# (C:/Users/Admin1/AppData/Roaming/PolarisOffice/ETemp/...). Clear it so
# the descr attribute is dropped entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Type -eq 13) {
        # msoPicture
        if ($shp.AlternativeText -like "*PolarisOffice*") {
            $shp.AlternativeText = ""
        }
    }
}
